$wb = $excel.ActiveWorkbook

# --- Sheet "Encabezado" (header / factura) ---
$wsEnc = $wb.Worksheets.Item("Encabezado")

# New recovered invoice header values (row 2)
$wsEnc.Range("A2").Value = "01135084"
$wsEnc.Range("B2").Value = "T010005360"
$wsEnc.Range("C2").Value = "0001"
$wsEnc.Range("D2").Value = "T3C1"
$wsEnc.Range("E2").Value = "MJIRON"
$wsEnc.Range("F2").Value = "CT3000000000012"
$wsEnc.Range("G2").Value = "T03"
$wsEnc.Range("H2").Value = 36.5635

# Column width clean-up (values line up with the engine's ColumnWidth -> stored
# width offset of 5/6, so subtracting it reproduces the exact target widths)
$wsEnc.Columns.Item(1).ColumnWidth = 15 - 0.8333333333333333
$wsEnc.Range("B1:C1").ColumnWidth = 18 - 0.8333333333333333
$wsEnc.Columns.Item(4).ColumnWidth = 17 - 0.8333333333333333
$wsEnc.Columns.Item(7).ColumnWidth = 12 - 0.8333333333333333
$wsEnc.Columns.Item(11).ColumnWidth = 22 - 0.8333333333333333

# --- Sheet "Detalles" (line items) ---
$wsDet = $wb.Worksheets.Item("Detalles")

# Remove the old second detail line (row 3) - only one line item remains
$wsDet.Rows(3).Delete()

# New recovered detail line (row 2)
$wsDet.Range("A2").Value = 0
$wsDet.Range("B2").Value = "01000019"
$wsDet.Range("C2").Value = "1.00"
$wsDet.Range("D2").Value = 0
$wsDet.Range("E2").Value = "ARROZ PAELLA VIGO 19OZ 012670 S/M 013052"
